# Actualización automática desde tarea programada
# Appends the new sensor reading row (row 9) to Sheet1, mirroring the
# style of the preceding rows (date/time cell keeps the custom date
# number format already applied to column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The scheduled task's export re-serialised the whole sheet, which
# re-quantised the existing row 8 timestamp by a sub-millisecond amount
# (floating point noise from the data source). Re-apply it verbatim.
$ws.Cells.Item(8, 1).Value = 45866.37525954861

$row = 9

$ws.Cells.Item($row, 1).Value = 45866.41686358196
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 2025
$ws.Cells.Item($row, 3).Value = 31
$ws.Cells.Item($row, 4).Value = 16.78
$ws.Cells.Item($row, 5).Value = 82.84999999999999
$ws.Cells.Item($row, 6).Value = 479.88
$ws.Cells.Item($row, 7).Value = 10.77
$ws.Cells.Item($row, 8).Value = "ESE"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = "10:00:17"
